$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force affected cells to Text format so numeric-looking strings
# (e.g. "36.30", "1.00") are preserved verbatim instead of being
# normalized to numbers by Excel when the .Value is assigned.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.149.45"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.17%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.322.84"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.86%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.52"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.65%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.98"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.71%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.508"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.14%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.72%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.30"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +6.17%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.52%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.03%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.82"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.55%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.93"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.86%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.683.66"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.92%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.285.64"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.21%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.799"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.34%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.055.66"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.17%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.13"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.31%  "

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.11%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0913"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.64%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.74%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.69"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.72%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.17"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.97%  "

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.02%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.06%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.52"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.13%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.99"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.59%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.31"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.89%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.21"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.79%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.59%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.95"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +8.82%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.27%  "

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.08%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.87"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.09%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.69%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.24%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.83"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.96%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.00%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.79"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.76%  "

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.39%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.995.94"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.09%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0290"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.73%  "

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.53%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.20"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.81%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.60"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.63%  "

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.31%  "

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.25"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.93%  "

# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.32"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +8.82%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.548.95"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.92%  "

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.55"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.86%  "
